$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the Korean (ko_KR) localized text.
# Each entry below corresponds to one <si> text edited in the diff.

$ws.Cells.Item(33, 4).Value = "[name=""웨이 옌우""]  정찰팀이 보고해주더군, 용문을 떠난 감염자 일부가 코어 쪽으로 이동 중이라고.
"

$ws.Cells.Item(71, 4).Value = "[name=""웨이 옌우""]  라이타니아든 사르곤이든 사미든…… 어느 곳이든 당신을 숨겨놓겠어.
"

$ws.Cells.Item(78, 4).Value = "[name=""후미즈키""]  웨이 옌우, 누누히 말했잖아요…… 전 괜찮다고.
"

$ws.Cells.Item(82, 4).Value = "[name=""후미즈키""]  하지만 아시잖아요, 당신은 절 막을 수 없다는 걸. 지금까지 당신은, 단 한 번도 저를 막아낸 적이 없었죠.
"

$ws.Cells.Item(106, 4).Value = "[name=""후미즈키""]  ……제가 하는 모든 일이 전부 당신만을 위한 건 아니라구요.
"

$ws.Cells.Item(124, 4).Value = "[name=""후미즈키""]  ……당신…… 이러지 말아요!
"
